# Update cryptocurrency price/volume snapshot values (Sheet1, columns D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.152.62'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').Value = '2.272.35'
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '298.31'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -2.62%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '94.67'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -5.70%  '
$ws.Range('E7').Value = '  +0.05%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.493'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -3.65%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.490'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -3.75%  '
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('E11').Value = '  -1.07%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '48.06'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -8.48%  '
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('E14').Value = '  -3.04%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '15.68'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '2.625.48'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').Value = '2.268.01'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('E18').Value = '  -4.59%  '
$ws.Range('D19').Value = '42.154.25'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('E20').Value = '  -2.44%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '11.37'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -3.46%  '
$ws.Range('E22').Value = '  -3.80%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '66.59'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -1.84%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '232.42'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -4.36%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '23.80'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -6.05%  '
$ws.Range('E29').Value = '  -1.16%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '166.13'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.48%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '33.69'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.83%  '
$ws.Range('E33').Value = '  -0.08%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '4.90'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -3.92%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.48'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('E36').Value = '  -5.51%  '
$ws.Range('E37').Value = '  -5.01%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '16.02'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -8.60%  '
$ws.Range('E39').Value = '  -4.64%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.0987'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -3.38%  '
$ws.Range('E41').Value = '  -3.65%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.72'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -7.12%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.40'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -7.73%  '
$ws.Range('D44').Value = '1.962.23'
$ws.Range('E45').Value = '  -2.93%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '17.45'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -6.13%  '
$ws.Range('E47').Value = '  -5.95%  '
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('D49').Value = '2.497.00'
$ws.Range('E49').Value = '  -2.22%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '51.96'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -6.95%  '
$ws.Range('E51').Value = '  -5.11%  '
